$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 51

$ws.Cells.Item($row, 1).Value = "IOLCNV"
$ws.Cells.Item($row, 2).Value = "Kit de engranaje de embrague para bandeja Epson"
$ws.Cells.Item($row, 3).Value = "L3100 L3101 L3110 L3115 L3116 L3150 L3151 L3156 L3158 L3160 L4150 L4160 L4165 L5190 L6178 L6198"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 5
$ws.Cells.Item($row, 7).Value = 4
$ws.Cells.Item($row, 8).Formula = "=(E51-D51)*G51"
$ws.Cells.Item($row, 9).Formula = "=D51*F51"
$ws.Cells.Item($row, 10).Value = 0
